$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Triplex_Fault"
$ws.Range("B7").Value = "Triplex Fault"
$ws.Range("C7").Value = "%M20"
$ws.Range("D7").Value = 20
$ws.Range("E7").Formula = "=D7+10001"
$ws.Range("F7").Value = "Boolean"
$ws.Range("H7").Value = "Alarm"
$ws.Range("I7").Value = "Ok"
$ws.Range("J7").Value = "Binary"
$ws.Range("K7").Value = "ro"
$ws.Range("L7").Value = "Pump 1"

$ws.Range("A8").Value = "South_Horizontal_Heat_Treat_Lvl_SW"
$ws.Range("B8").Value = "South Horizontal Heater Lvl SW"
$ws.Range("C8").Value = 1099
$ws.Range("E8").Value = 1099
$ws.Range("F8").Value = "Boolean"
$ws.Range("H8").Value = "Alarm"
$ws.Range("I8").Value = "Ok"
$ws.Range("J8").Value = "Binary"
$ws.Range("K8").Value = "RO"
$ws.Range("L8").Value = "South Heater Treater"

$ws.Range("A7").Select()
